$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E101").Value = "{'List[any]', 'empty'}"

$ws.Range("E102").Value = "List[any]"
$ws.Range("F102").Value = "Neutral"
$ws.Range("F102").Interior.Color = 42495

$ws.Range("D125").Value = 10

$ws.Range("C126").Value = ""
$ws.Range("D126").Value = ""
$ws.Range("E126").Value = "Scalpel Accuracy:"
$ws.Range("F126").Value = 91.86999999999999

$ws.Range("E127").Value = "Accuracy vs PyType"
$ws.Range("F127").Value = 240
